# commit: "handle dividion by 0 and output step with no previous steps"
#
# The test log sheet gets two new test-case rows (17: a Number/Calculus
# combo exercising the divide-by-zero fix, 18: more invalid-input /
# flow-analysis coverage), a label fix in the existing "Output" step
# (row 20, column D: "4. Output..." -> "3. Output..." since it now has
# no previous steps), and a relabeling of test 10's steps (row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test 17 (row 21): Number -> Calculus(Step 1) -> Number -> Output ---
$ws.Range("B21").Value = "1. Number(description1)(11)(1)"
$ws.Range("C21").Value = "2. Calculus()(Step 1)(1)"

# Existing "Output" step on row 20 now reads as step 3 (no previous steps
# to reference), and test 17 continues with steps 3 and 4.
$ws.Range("D20").Value = "3. Output(all previous steps)"
$ws.Range("D21").Value = "3. Number(description2)(22)(1)"
$ws.Range("E21").Value = "4. Output()(Steps 1, 2 and 3)"

# --- Test 18 (row 22): divide-by-zero handling + bad input coverage ---
$ws.Range("B22").Value = "1. TextFile(description1)(test.txt)(1)"
$ws.Range("D22").Value = "3. Number(description3)(11)(1)"
$ws.Range("C22").Value = "2. CsvFile(description2)(test.csv)(0)"
$ws.Range("E22").Value = "4. Entered bad input a couple of times"
$ws.Range("F22").Value = "5. Ran flow analysis(details) option"

# --- Test 10 (row 14): relabeled Number steps (now covers 0-division) ---
$ws.Range("B14").Value = "1. Number(description1)(111)(1)"
$ws.Range("C14").Value = "3. Number(description2)(0)(0)"

# Columns E/F grew to fit the new, longer step descriptions.
$ws.Columns.Item(5).ColumnWidth = 32.166666666666664
$ws.Columns.Item(6).ColumnWidth = 29.666666666666668

$ws.Range("C15").Select()
